$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "ENGLISH": insert one new vocabulary row in the middle (row 68, the
# word "extricate"), then append nine more new vocabulary rows at the end
# (rows 136-144).
# ---------------------------------------------------------------------------
$wsEnglish = $wb.Worksheets.Item("ENGLISH")

# Insert a new blank row at position 68 - this shifts every row that used to
# be at 68.. down by one (old row 68 "slumber" becomes row 69, etc).
$wsEnglish.Rows.Item(68).Insert()

$wsEnglish.Range("A68").Value = "extricate"
$wsEnglish.Range("B68").Value = "free from a constraint or difficulty"
$wsEnglish.Range("C68").Value = ""
$wsEnglish.Range("D68").Value = ""
$wsEnglish.Range("E68").Value = 0
$wsEnglish.Range("F68").Value = "2021-11-20 14:57:27.485397"
$wsEnglish.Range("G68").Value = ""

$wsEnglish.Range("A136").Value = "restrained"
$wsEnglish.Range("B136").Value = ""
$wsEnglish.Range("C136").Value = "self-controlled"
$wsEnglish.Range("D136").Value = ""
$wsEnglish.Range("E136").Value = 0
$wsEnglish.Range("F136").Value = "2021-11-20 15:03:49.858279"
$wsEnglish.Range("G136").Value = ""

$wsEnglish.Range("A137").Value = "condense"
$wsEnglish.Range("B137").Value = "make denser or more concentrated"
$wsEnglish.Range("C137").Value = ""
$wsEnglish.Range("D137").Value = ""
$wsEnglish.Range("E137").Value = 0
$wsEnglish.Range("F137").Value = "2021-11-20 15:04:42.212772"
$wsEnglish.Range("G137").Value = ""

$wsEnglish.Range("A138").Value = "apt"
$wsEnglish.Range("B138").Value = ""
$wsEnglish.Range("C138").Value = "inclined;suitable"
$wsEnglish.Range("D138").Value = ""
$wsEnglish.Range("E138").Value = 0
$wsEnglish.Range("F138").Value = "2021-11-20 15:05:21.928371"
$wsEnglish.Range("G138").Value = ""

$wsEnglish.Range("A139").Value = "sentiment"
$wsEnglish.Range("B139").Value = "a view or opinion that is held or expressed"
$wsEnglish.Range("C139").Value = "view;feeling"
$wsEnglish.Range("D139").Value = ""
$wsEnglish.Range("E139").Value = 0
$wsEnglish.Range("F139").Value = "2021-11-20 15:06:18.444516"
$wsEnglish.Range("G139").Value = ""

$wsEnglish.Range("A140").Value = "fuming"
$wsEnglish.Range("B140").Value = "feeling, showing or expressing great anger"
$wsEnglish.Range("C140").Value = ""
$wsEnglish.Range("D140").Value = ""
$wsEnglish.Range("E140").Value = 0
$wsEnglish.Range("F140").Value = "2021-11-20 15:07:39.355342"
$wsEnglish.Range("G140").Value = ""

$wsEnglish.Range("A141").Value = "insidiously"
$wsEnglish.Range("B141").Value = "in a gradual, subtle way, but with harmful effects"
$wsEnglish.Range("C141").Value = ""
$wsEnglish.Range("D141").Value = ""
$wsEnglish.Range("E141").Value = 0
$wsEnglish.Range("F141").Value = "2021-11-20 15:08:28.761904"
$wsEnglish.Range("G141").Value = ""

$wsEnglish.Range("A142").Value = "fret"
$wsEnglish.Range("B142").Value = "be constantly or visibly anxious"
$wsEnglish.Range("C142").Value = "worry;trouble"
$wsEnglish.Range("D142").Value = ""
$wsEnglish.Range("E142").Value = 0
$wsEnglish.Range("F142").Value = "2021-11-20 15:10:17.940034"
$wsEnglish.Range("G142").Value = ""

$wsEnglish.Range("A143").Value = "unsolicited"
$wsEnglish.Range("B143").Value = "not asked for; given or done voluntarily"
$wsEnglish.Range("C143").Value = "uninvited"
$wsEnglish.Range("D143").Value = ""
$wsEnglish.Range("E143").Value = 0
$wsEnglish.Range("F143").Value = "2021-11-20 15:12:00.23661"
$wsEnglish.Range("G143").Value = ""

$wsEnglish.Range("A144").Value = "resemble"
$wsEnglish.Range("B144").Value = ""
$wsEnglish.Range("C144").Value = "look like"
$wsEnglish.Range("D144").Value = ""
$wsEnglish.Range("E144").Value = 0
$wsEnglish.Range("F144").Value = "2021-11-20 15:12:28.779915"
$wsEnglish.Range("G144").Value = ""

# ---------------------------------------------------------------------------
# Sheet "NOTES": append seven new quote/content rows at the end (rows 29-35).
# ---------------------------------------------------------------------------
$wsNotes = $wb.Worksheets.Item("NOTES")

$wsNotes.Range("A29").Value = "Get over the fear of waste"
$wsNotes.Range("B29").Value = "essentialism"

$wsNotes.Range("A30").Value = "Stop making casual commitments"
$wsNotes.Range("B30").Value = "essentialism"

$wsNotes.Range("A31").Value = "Pause before you speak"
$wsNotes.Range("B31").Value = ""

$wsNotes.Range("A32").Value = "get over the fear of missing out"
$wsNotes.Range("B32").Value = ""

$wsNotes.Range("A33").Value = "I saw the angel in the marble and carved until I set him free"
$wsNotes.Range("B33").Value = ""

$wsNotes.Range("A34").Value = "No is a complete sentence"
$wsNotes.Range("B34").Value = ""

$wsNotes.Range("A35").Value = "If you don''t set boundaries - there won''t be any"
$wsNotes.Range("B35").Value = ""
